# ---------------------------------------------------------------------------
# Change Navigation bar style, logo and Add contact, selfstudy.html template
#
# 1. Re-cache the footer "date" placeholder (datetimeFigureOut field) on the
#    slide master and all 11 slide layouts: "2016. 9. 14." -> "2016-11-29".
# 2. Add a new logo/navbar group (4 text boxes) onto the last existing slide
#    (slide 4, currently blank).
# 3. Append a brand-new blank slide (slide 5).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Step 1: update the cached date field text wherever it is a "date"
# placeholder (ppPlaceholderDate = 16), on the slide master and every layout.
# ---------------------------------------------------------------------------
$newDate = "2016-11-29"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = -1 }
        if ($phType -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------------
# Step 2: build the logo / navbar group on slide 4 (last slide, blank).
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Two throw-away textboxes so the shape-id counter lands on 4/5/6 and then
# 8 for the real shapes (mirrors the authoring history captured in the
# target deck, where ids 2,3 and 7 belong to shapes that no longer exist).
$skip1 = $s4.Shapes.AddTextbox(1, 0, 0, 1, 1)
$skip2 = $s4.Shapes.AddTextbox(1, 0, 0, 1, 1)
$skip1.Delete()
$skip2.Delete()

$EMU = 12700.0

# --- TextBox 3 (id 4): "YOUNG / HUN / JUNG" -------------------------------
$tb3 = $s4.Shapes.AddTextbox(1, 4998198/$EMU, 3036880/$EMU, 1796068/$EMU, 770980/$EMU)
$tb3.Fill.Visible = $false
$tb3.TextFrame.WordWrap = -1
$tb3.TextFrame.AutoSize = 1
$tb3.TextFrame.TextRange.Text = "YOUNG`rHUN`rJUNG"
$tb3.TextFrame.TextRange.Font.Bold = $true
$tb3.TextFrame.TextRange.Font.NameAscii = "Apple SD Gothic Neo"
$tb3.TextFrame.TextRange.Font.NameFarEast = "Apple SD Gothic Neo"
$tb3.TextFrame.TextRange.Font.NameComplexScript = "Apple SD Gothic Neo"
for ($i = 1; $i -le 3; $i++) {
    $pf = $tb3.TextFrame.TextRange.Paragraphs($i, 1).ParagraphFormat
    $pf.LineRuleWithin = $true
    $pf.SpaceWithin = 0.8
}

# --- TextBox 4 (id 5): "JUNG / YOUNG / HUN" (right aligned, orange) -------
$tb4 = $s4.Shapes.AddTextbox(1, 5723467/$EMU, 3054301/$EMU, 1568732/$EMU, 770980/$EMU)
$tb4.Fill.Visible = $false
$tb4.TextFrame.WordWrap = -1
$tb4.TextFrame.AutoSize = 1
$tb4.TextFrame.TextRange.Text = "JUNG`rYOUNG`rHUN"
$tb4.TextFrame.TextRange.Font.Bold = $true
$tb4.TextFrame.TextRange.Font.NameAscii = "Apple SD Gothic Neo"
$tb4.TextFrame.TextRange.Font.NameFarEast = "Apple SD Gothic Neo"
$tb4.TextFrame.TextRange.Font.NameComplexScript = "Apple SD Gothic Neo"
$tb4.TextFrame.TextRange.Font.Color.RGB = 18175  # FF4600
for ($i = 1; $i -le 3; $i++) {
    $para = $tb4.TextFrame.TextRange.Paragraphs($i, 1)
    $para.ParagraphFormat.Alignment = 3
    $pf = $para.ParagraphFormat
    $pf.LineRuleWithin = $true
    $pf.SpaceWithin = 0.8
}

# --- TextBox 5 (id 6): top decorative "l l l" bar -------------------------
$tb5 = $s4.Shapes.AddTextbox(1, 4901940/$EMU, 2570600/$EMU, 2422688/$EMU, 230832/$EMU)
$tb5.Fill.Visible = $false
$tb5.TextFrame.WordWrap = -1
$tb5.TextFrame.AutoSize = 1
$tf5 = $tb5.TextFrame
$tf5.TextRange.Text = "l l l l "
$tf5.TextRange.InsertAfter("l l l l l ") | Out-Null
$tf5.TextRange.InsertAfter("l l l l l l l") | Out-Null
$tf5.TextRange.InsertAfter(" ") | Out-Null
$tf5.TextRange.InsertAfter("l l l l l l l l l l ") | Out-Null
$tf5.TextRange.InsertAfter("l l l l l l l ") | Out-Null
$tf5.TextRange.Font.Size = 9
$tf5.TextRange.Characters(1, 8).Font.Color.RGB = 3881983    # FF3B3B
$tf5.TextRange.Characters(9, 10).Font.Color.RGB = 180       # B40000
$tf5.TextRange.Characters(19, 13).Font.Color.RGB = 201111   # 971103
$tf5.TextRange.Characters(32, 1).Font.Color.RGB = 158       # 9E0000
$tf5.TextRange.Characters(33, 20).Font.Color.RGB = 100      # 640000

# --- TextBox 7 (id 8): bottom decorative "l l l" bar (rotated 180deg) -----
$skip3 = $s4.Shapes.AddTextbox(1, 0, 0, 1, 1)
$skip3.Delete()

$tb7 = $s4.Shapes.AddTextbox(1, 4901940/$EMU, 3945313/$EMU, 2422688/$EMU, 230832/$EMU)
$tb7.Fill.Visible = $false
$tb7.TextFrame.WordWrap = -1
$tb7.TextFrame.AutoSize = 1
$tb7.Rotation = 180
$tf7 = $tb7.TextFrame
$tf7.TextRange.Text = "l l l l "
$tf7.TextRange.InsertAfter("l l l l l ") | Out-Null
$tf7.TextRange.InsertAfter("l l l l l l l") | Out-Null
$tf7.TextRange.InsertAfter(" ") | Out-Null
$tf7.TextRange.InsertAfter("l l l l l l l l l l ") | Out-Null
$tf7.TextRange.InsertAfter("l l l l l l l ") | Out-Null
$tf7.TextRange.Font.Size = 9
$tf7.TextRange.Characters(1, 8).Font.Color.RGB = 3881983    # FF3B3B
$tf7.TextRange.Characters(9, 10).Font.Color.RGB = 180       # B40000
$tf7.TextRange.Characters(19, 13).Font.Color.RGB = 201111   # 971103
$tf7.TextRange.Characters(32, 1).Font.Color.RGB = 158       # 9E0000
$tf7.TextRange.Characters(33, 20).Font.Color.RGB = 100      # 640000

# Group the four shapes together and name the group like the authored deck.
$logoGroup = $s4.Shapes.Range(@($tb3.Name, $tb4.Name, $tb5.Name, $tb7.Name)).Group()
$logoGroup.Name = "그룹 8"

# ---------------------------------------------------------------------------
# Step 3: add the new, blank 5th slide.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 12)

Write-Output ("Slides: " + $p.Slides.Count)
